$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-07 Saturday" "2024-12-08 Sunday"

Replace-Text "421÷3=140, 1" "816÷9=90, 6"
Replace-Text "267÷5=53, 2" "382÷3=127, 1"
Replace-Text "643÷9=71, 4" "144÷7=20, 4"
Replace-Text "330÷9=36, 6" "562÷6=93, 4"
Replace-Text "190÷3=63, 1" "156÷8=19, 4"

Replace-Text "324÷3=108, 0" "481÷4=120, 1"
Replace-Text "607÷4=151, 3" "403÷5=80, 3"
Replace-Text "869÷2=434, 1" "766÷4=191, 2"
Replace-Text "648÷2=324, 0" "600÷3=200, 0"
Replace-Text "623÷3=207, 2" "560÷9=62, 2"

Replace-Text "946÷9=105, 1" "715÷9=79, 4"
Replace-Text "845÷4=211, 1" "438÷5=87, 3"
Replace-Text "711÷2=355, 1" "935÷7=133, 4"
Replace-Text "365÷7=52, 1" "901÷8=112, 5"
Replace-Text "437÷3=145, 2" "417÷3=139, 0"

Replace-Text "409÷7=58, 3" "613÷4=153, 1"
Replace-Text "320÷3=106, 2" "711÷6=118, 3"
Replace-Text "359÷7=51, 2" "982÷7=140, 2"
Replace-Text "903÷3=301, 0" "673÷2=336, 1"
Replace-Text "479÷2=239, 1" "618÷8=77, 2"

Replace-Text "446÷2=223, 0" "749÷9=83, 2"
Replace-Text "862÷7=123, 1" "438÷2=219, 0"
Replace-Text "436÷5=87, 1" "471÷4=117, 3"
Replace-Text "432÷2=216, 0" "366÷3=122, 0"
Replace-Text "579÷8=72, 3" "674÷7=96, 2"

Write-Output "Done"
